$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text-formatted
# representation (e.g. "1.00", "66.018.03") instead of Excel auto-
# converting numeric-looking strings into real numbers.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value2 = '66.018.03'
$ws.Range('E2').Value2 = '  -1.36%  '
$ws.Range('D3').Value2 = '3.476.93'
$ws.Range('E3').Value2 = '  -0.10%  '
$ws.Range('D4').Value2 = '1.00'
$ws.Range('E4').Value2 = '  +0.08%  '
$ws.Range('D5').Value2 = '603.16'
$ws.Range('E5').Value2 = '  +0.39%  '
$ws.Range('D6').Value2 = '142.84'
$ws.Range('E6').Value2 = '  -3.21%  '
$ws.Range('D7').Value2 = '3.473.87'
$ws.Range('E7').Value2 = '  -0.12%  '
$ws.Range('E8').Value2 = '  -0.09%  '
$ws.Range('D9').Value2 = '0.474'
$ws.Range('E9').Value2 = '  -0.97%  '
$ws.Range('D10').Value2 = '8.21'
$ws.Range('E10').Value2 = '  +7.41%  '
$ws.Range('E11').Value2 = '  -4.80%  '
$ws.Range('D12').Value2 = '0.411'
$ws.Range('E12').Value2 = '  -2.54%  '
$ws.Range('D13').Value2 = '4.070.22'
$ws.Range('E13').Value2 = '  +0.07%  '
$ws.Range('D14').Value2 = '0.0000203'
$ws.Range('E14').Value2 = '  -4.16%  '
$ws.Range('D15').Value2 = '30.31'
$ws.Range('E15').Value2 = '  -2.82%  '
$ws.Range('D16').Value2 = '3.475.49'
$ws.Range('E16').Value2 = '  -0.02%  '
$ws.Range('D17').Value2 = '0.116'
$ws.Range('D18').Value2 = '66.081.61'
$ws.Range('E18').Value2 = '  -1.17%  '
$ws.Range('D19').Value2 = '10.40'
$ws.Range('E19').Value2 = '  +3.30%  '
$ws.Range('D20').Value2 = '6.15'
$ws.Range('E20').Value2 = '  -3.67%  '
$ws.Range('D21').Value2 = '14.70'
$ws.Range('E21').Value2 = '  -3.48%  '
$ws.Range('D22').Value2 = '419.80'
$ws.Range('E22').Value2 = '  -3.11%  '
$ws.Range('D23').Value2 = '0.586'
$ws.Range('E23').Value2 = '  -3.03%  '
$ws.Range('D24').Value2 = '77.46'
$ws.Range('E24').Value2 = '  -2.03%  '
$ws.Range('B25').Value2 = 'WrappedeETH'
$ws.Range('C25').Value2 = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value2 = '3.627.27'
$ws.Range('E25').Value2 = '  +0.30%  '
$ws.Range('B26').Value2 = 'Dai'
$ws.Range('C26').Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value2 = '1.00'
$ws.Range('E26').Value2 = '  -0.02%  '
$ws.Range('B27').Value2 = 'PEPE'
$ws.Range('C27').Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value2 = '0.0000115'
$ws.Range('E27').Value2 = '  -3.35%  '
$ws.Range('B28').Value2 = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value2 = '9.39'
$ws.Range('E28').Value2 = '  -4.13%  '
$ws.Range('B29').Value2 = 'RenderToken'
$ws.Range('C29').Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value2 = '7.94'
$ws.Range('E29').Value2 = '  -4.97%  '
$ws.Range('B30').Value2 = 'PancakeSwap'
$ws.Range('C30').Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value2 = '2.46'
$ws.Range('E30').Value2 = '  -0.55%  '
$ws.Range('B31').Value2 = 'Binance-PegBSC-USD'
$ws.Range('C31').Value2 = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value2 = '1.00'
$ws.Range('E31').Value2 = '  +0.10%  '
$ws.Range('B32').Value2 = 'Kaspa'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value2 = '0.162'
$ws.Range('E32').Value2 = '  -3.10%  '
$ws.Range('B33').Value2 = 'Fetch.AI'
$ws.Range('C33').Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value2 = '1.48'
$ws.Range('E33').Value2 = '  -6.10%  '
$ws.Range('B34').Value2 = 'EthereumClassic'
$ws.Range('C34').Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value2 = '25.06'
$ws.Range('E34').Value2 = '  -0.89%  '
$ws.Range('B35').Value2 = 'RenzoRestakedETH'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D35').Value2 = '3.474.48'
$ws.Range('E35').Value2 = '  +0.20%  '
$ws.Range('B36').Value2 = 'USDe'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D36').Value2 = '1.00'
$ws.Range('E36').Value2 = '  -0.08%  '
$ws.Range('B37').Value2 = 'ImmutableX'
$ws.Range('C37').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value2 = '1.70'
$ws.Range('E37').Value2 = '  -5.12%  '
$ws.Range('B38').Value2 = 'NEARProtocol'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value2 = '5.54'
$ws.Range('E38').Value2 = '  -6.07%  '
$ws.Range('B39').Value2 = 'Aptos'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value2 = '7.62'
$ws.Range('E39').Value2 = '  -3.27%  '
$ws.Range('B40').Value2 = 'FirstDigitalUSD'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value2 = '0.999'
$ws.Range('E40').Value2 = '  +0.05%  '
$ws.Range('B41').Value2 = 'Monero'
$ws.Range('C41').Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value2 = '169.76'
$ws.Range('E41').Value2 = '  -2.30%  '
$ws.Range('B42').Value2 = 'Hedera'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').Value2 = '0.0866'
$ws.Range('E42').Value2 = '  -1.90%  '
$ws.Range('B43').Value2 = 'Mantle'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value2 = '0.889'
$ws.Range('E43').Value2 = '  -0.77%  '
$ws.Range('B44').Value2 = 'Filecoin'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value2 = '5.10'
$ws.Range('E44').Value2 = '  -5.17%  '
$ws.Range('B45').Value2 = 'Stacks'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value2 = '1.91'
$ws.Range('E45').Value2 = '  -8.44%  '
$ws.Range('B46').Value2 = 'OKB'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value2 = '45.66'
$ws.Range('E46').Value2 = '  -1.65%  '
$ws.Range('B47').Value2 = 'InjectiveProtocol'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value2 = '26.12'
$ws.Range('E47').Value2 = '  -9.25%  '
$ws.Range('B48').Value2 = 'ONDO'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value2 = '1.20'
$ws.Range('E48').Value2 = '  -3.50%  '
$ws.Range('B49').Value2 = 'Cosmos'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value2 = '7.12'
$ws.Range('E49').Value2 = '  -4.31%  '
$ws.Range('B50').Value2 = 'dogwifhat'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value2 = '2.34'
$ws.Range('E50').Value2 = '  -2.71%  '
$ws.Range('B51').Value2 = 'SuiNetwork'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').Value2 = '0.930'
$ws.Range('E51').Value2 = '  -4.18%  '
